$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D cells whose new value would otherwise be auto-converted by Excel
# into a numeric value (e.g. "215.57" -> 215.57). Force them to text first so
# the stored value / cell type matches the original (inline string) content,
# then restore the default "Normal" style so no stray number format lingers.
$numericLookingAddrs = @("D5", "D6", "D14", "D15", "D16", "D18", "D19", "D23", "D24", "D25", "D27", "D30", "D36", "D38", "D40", "D42", "D47", "D49", "D50", "D51")
$numericLookingValues = @{
    "D5" = "215.57"
    "D6" = "0.517"
    "D14" = "4.11"
    "D15" = "0.534"
    "D16" = "66.38"
    "D18" = "8.19"
    "D19" = "236.07"
    "D23" = "9.25"
    "D24" = "2.11"
    "D25" = "147.47"
    "D27" = "16.52"
    "D30" = "0.0498"
    "D36" = "2.38"
    "D38" = "0.915"
    "D40" = "1.05"
    "D42" = "67.88"
    "D47" = "90.36"
    "D49" = "0.103"
    "D50" = "8.03"
    "D51" = "0.0507"
}
foreach ($addr in $numericLookingAddrs) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericLookingValues[$addr]
    $cell.Style = "Normal"
}

# Remaining D/E cell updates (these stay text automatically, no special handling needed)
$ws.Range("D2").Value = "27.047.33"
$ws.Range("E2").Value = "  +0.56%  "
$ws.Range("D3").Value = "1.677.78"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E5").Value = "  +0.13%  "
$ws.Range("E6").Value = "  -1.03%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +2.41%  "
$ws.Range("E9").Value = "  +5.40%  "
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("D12").Value = "1.915.39"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "1.700.10"
$ws.Range("E13").Value = "  +1.80%  "
$ws.Range("E14").Value = "  +0.85%  "
$ws.Range("E15").Value = "  +1.45%  "
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "27.044.99"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("E20").Value = "  +0.45%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("E23").Value = "  +1.29%  "
$ws.Range("E24").Value = "  -2.59%  "
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("E26").Value = "  +2.07%  "
$ws.Range("E27").Value = "  +4.09%  "
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  +0.19%  "
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("E32").Value = "  +0.98%  "
$ws.Range("D33").Value = "1.540.84"
$ws.Range("E33").Value = "  +6.45%  "
$ws.Range("E34").Value = "  +1.17%  "
$ws.Range("E35").Value = "  +4.98%  "
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("E37").Value = "  +0.90%  "
$ws.Range("E38").Value = "  +1.25%  "
$ws.Range("E39").Value = "  +2.19%  "
$ws.Range("E40").Value = "  +4.74%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +2.89%  "
$ws.Range("E43").Value = "  -3.58%  "
$ws.Range("E44").Value = "  -1.89%  "
$ws.Range("D45").Value = "1.821.80"
$ws.Range("E45").Value = "  +0.77%  "
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("E48").Value = "  +0.12%  "
$ws.Range("E49").Value = "  +1.71%  "
$ws.Range("E50").Value = "  +6.19%  "
$ws.Range("E51").Value = "  +0.17%  "
